# Refresh the stochastic scenario-comparison numbers (Raw_Data + Summary)
# produced by a re-run of the forest-growth simulation integration test.
$wb = $excel.ActiveWorkbook

$rawData = $wb.Worksheets.Item("Raw_Data")
$summary = $wb.Worksheets.Item("Summary")

# --- Raw_Data sheet: per-period simulation output for each scenario ---
$rawData.Range("B2").Value = 0.7202346731898849
$rawData.Range("C2").Value = 0.513911475849385
$rawData.Range("D2").Value = 0.5044358179472358
$rawData.Range("G2").Value = 0.05739913488953671
$rawData.Range("H2").Value = 4.265262369620648
$rawData.Range("L2").Value = 93.62013125626211
$rawData.Range("A3").Value = 491
$rawData.Range("B3").Value = 22.20308141284221
$rawData.Range("C3").Value = 2.879401224404767
$rawData.Range("D3").Value = 2.879374177297288
$rawData.Range("E3").Value = 22.63580122526506
$rawData.Range("F3").Value = 22.56889308861816
$rawData.Range("G3").Value = 1.769475578473134
$rawData.Range("H3").Value = 66.56738536087609
$rawData.Range("J3").Value = 1.386820528351585
$rawData.Range("L3").Value = 321.3059620113886
$rawData.Range("A4").Value = 480
$rawData.Range("B4").Value = 53.47293081480863
$rawData.Range("C4").Value = 4.5194197665066
$rawData.Range("D4").Value = 4.519151007196385
$rawData.Range("E4").Value = 36.09277117852089
$rawData.Range("F4").Value = 36.0284581218981
$rawData.Range("G4").Value = 2.946718756668255
$rawData.Range("H4").Value = 134.1679981583152
$rawData.Range("J4").Value = 2.795166628298233
$rawData.Range("L4").Value = 972.6668686620575
$rawData.Range("A5").Value = 467
$rawData.Range("B5").Value = 95.00177720403157
$rawData.Range("C5").Value = 6.107222321845692
$rawData.Range("D5").Value = 6.10554479002972
$rawData.Range("E5").Value = 48.31721644460434
$rawData.Range("F5").Value = 48.25153233423276
$rawData.Range("G5").Value = 4.202724579741611
$rawData.Range("H5").Value = 211.6391785164422
$rawData.Range("J5").Value = 4.409149552425879
$rawData.Range("L5").Value = 2188.437981188561
$rawData.Range("M5").Value = 614.0344210662722
$rawData.Range("A6").Value = 460
$rawData.Range("B6").Value = 137.8001182480799
$rawData.Range("C6").Value = 7.411091660251636
$rawData.Range("D6").Value = 7.405522824688475
$rawData.Range("E6").Value = 58.02817591230716
$rawData.Range("F6").Value = 57.95892898370582
$rawData.Range("G6").Value = 5.407108944112923
$rawData.Range("H6").Value = 284.3928922783617
$rawData.Range("J6").Value = 5.924851922465868
$rawData.Range("L6").Value = 3747.081408254661
$rawData.Range("M6").Value = 1478.616647725951
$rawData.Range("B7").Value = 1.151402426919578
$rawData.Range("C7").Value = 0.5136944600970971
$rawData.Range("D7").Value = 0.5037685693587338
$rawData.Range("G7").Value = 0.0917610685447697
$rawData.Range("H7").Value = 6.819795041725524
$rawData.Range("L7").Value = 149.7917639966689
$rawData.Range("A8").Value = 777
$rawData.Range("B8").Value = 35.11561962361235
$rawData.Range("C8").Value = 2.878564492392675
$rawData.Range("D8").Value = 2.878531130556662
$rawData.Range("E8").Value = 22.63111160174193
$rawData.Range("F8").Value = 22.56248150356302
$rawData.Range("G8").Value = 2.798540895814356
$rawData.Range("H8").Value = 105.2927433121898
$rawData.Range("J8").Value = 2.193598819003953
$rawData.Range("L8").Value = 508.1504159660739
$rawData.Range("A9").Value = 763
$rawData.Range("B9").Value = 82.11858239032379
$rawData.Range("C9").Value = 4.442165738462506
$rawData.Range("D9").Value = 4.441632988807243
$rawData.Range("E9").Value = 36.02203883103742
$rawData.Range("F9").Value = 35.95680124014695
$rawData.Range("G9").Value = 4.588322905533941
$rawData.Range("H9").Value = 207.4503275129931
$rawData.Range("J9").Value = 4.321881823187356
$rawData.Range("L9").Value = 1495.890735317712
$rawData.Range("A10").Value = 751
$rawData.Range("B10").Value = 139.3622383003088
$rawData.Range("C10").Value = 5.832958591087717
$rawData.Range("D10").Value = 5.829649439338227
$rawData.Range("E10").Value = 48.11454877754961
$rawData.Range("F10").Value = 48.04658287481557
$rawData.Range("G10").Value = 6.357301543717624
$rawData.Range("H10").Value = 316.1488532056907
$rawData.Range("J10").Value = 6.586434441785224
$rawData.Range("L10").Value = 3209.655512182477
$rawData.Range("M10").Value = 810.1570025813929
$rawData.Range("A11").Value = 736
$rawData.Range("B11").Value = 190.303266909769
$rawData.Range("C11").Value = 6.885262781597936
$rawData.Range("D11").Value = 6.875513654645529
$rawData.Range("E11").Value = 57.7301768256572
$rawData.Range("F11").Value = 57.6570267755023
$rawData.Range("G11").Value = 7.798246157536356
$rawData.Range("H11").Value = 404.3339878014114
$rawData.Range("J11").Value = 8.423624745862737
$rawData.Range("L11").Value = 5167.488790218209
$rawData.Range("M11").Value = 1839.928815155501
$rawData.Range("B12").Value = 0.7137702800874285
$rawData.Range("C12").Value = 0.5115999974410029
$rawData.Range("D12").Value = 0.501725070229983
$rawData.Range("G12").Value = 0.05688395478855179
$rawData.Range("H12").Value = 4.234513393189103
$rawData.Range("L12").Value = 93.61716819672704
$rawData.Range("A13").Value = 492
$rawData.Range("B13").Value = 13.76953395402863
$rawData.Range("C13").Value = 2.265235580889897
$rawData.Range("D13").Value = 2.265212596266086
$rawData.Range("E13").Value = 17.95879625128244
$rawData.Range("F13").Value = 17.90496140681934
$rawData.Range("G13").Value = 1.097363631901914
$rawData.Range("H13").Value = 45.38601674753248
$rawData.Range("L13").Value = 204.8087256083742
$rawData.Range("A14").Value = 475
$rawData.Range("B14").Value = 37.13373621317103
$rawData.Range("C14").Value = 3.785940850933355
$rawData.Range("D14").Value = 3.785866775246867
$rawData.Range("E14").Value = 29.13301342225243
$rawData.Range("F14").Value = 29.08828149323029
$rawData.Range("G14").Value = 2.375764370630197
$rawData.Range("H14").Value = 99.92234049449684
$rawData.Range("J14").Value = 2.081715426968684
$rawData.Range("L14").Value = 583.7458474443873
$rawData.Range("A15").Value = 461
$rawData.Range("B15").Value = 71.89404937001953
$rawData.Range("C15").Value = 5.347272486405408
$rawData.Range("D15").Value = 5.346508437128667
$rawData.Range("E15").Value = 38.7493746968496
$rawData.Range("F15").Value = 38.70334980046432
$rawData.Range("G15").Value = 3.486282886248371
$rawData.Range("H15").Value = 168.7925567775223
$rawData.Range("J15").Value = 3.516511599531715
$rawData.Range("L15").Value = 1361.500772819123
$rawData.Range("M15").Value = 255.5878685663937
$rawData.Range("A16").Value = 452
$rawData.Range("B16").Value = 111.012337533814
$rawData.Range("C16").Value = 6.710466485806138
$rawData.Range("D16").Value = 6.707368506691897
$rawData.Range("E16").Value = 46.3931922815417
$rawData.Range("F16").Value = 46.34522099789339
$rawData.Range("G16").Value = 4.623335021611265
$rawData.Range("H16").Value = 238.2739350045039
$rawData.Range("J16").Value = 4.964040312593832
$rawData.Range("L16").Value = 2442.758324069539
$rawData.Range("M16").Value = 815.5289352177025
$rawData.Range("B17").Value = 0.687929221410262
$rawData.Range("C17").Value = 0.5022537429819555
$rawData.Range("D17").Value = 0.492733169036367
$rawData.Range("G17").Value = 0.04389351959216424
$rawData.Range("H17").Value = 4.111040084573148
$rawData.Range("L17").Value = 93.60532352792423
$rawData.Range("A18").Value = 495
$rawData.Range("B18").Value = 15.92202244767113
$rawData.Range("C18").Value = 2.428470441920808
$rawData.Range("D18").Value = 2.428446478652079
$rawData.Range("E18").Value = 22.63625906182645
$rawData.Range("F18").Value = 22.56831595571246
$rawData.Range("G18").Value = 1.01590916987219
$rawData.Range("H18").Value = 51.05806806878169
$rawData.Range("J18").Value = 1.04200138915881
$rawData.Range("L18").Value = 257.069849815854
$rawData.Range("A19").Value = 481
$rawData.Range("B19").Value = 42.83804835625742
$rawData.Range("C19").Value = 4.04090467497139
$rawData.Range("D19").Value = 4.040836549491504
$rawData.Range("E19").Value = 36.58707673776006
$rawData.Range("F19").Value = 36.5294401713959
$rawData.Range("G19").Value = 2.145886683648852
$rawData.Range("H19").Value = 112.3422811871288
$rawData.Range("J19").Value = 2.292699616063854
$rawData.Range("L19").Value = 807.0372324284087
$rawData.Range("A20").Value = 474
$rawData.Range("B20").Value = 79.25932786285384
$rawData.Range("C20").Value = 5.536971996804944
$rawData.Range("D20").Value = 5.536372678018728
$rawData.Range("E20").Value = 50.11560307827571
$rawData.Range("F20").Value = 50.0578465482267
$rawData.Range("G20").Value = 3.214844027538497
$rawData.Range("H20").Value = 183.5398729037619
$rawData.Range("J20").Value = 3.745711691913507
$rawData.Range("L20").Value = 1907.107672888115
$rawData.Range("M20").Value = 420.0997541497803
$rawData.Range("A21").Value = 469
$rawData.Range("B21").Value = 118.2477830309855
$rawData.Range("C21").Value = 6.799021078141265
$rawData.Range("D21").Value = 6.796704856080432
$rawData.Range("E21").Value = 63.47513064507812
$rawData.Range("F21").Value = 63.41420701930006
$rawData.Range("G21").Value = 4.275609357765627
$rawData.Range("H21").Value = 252.4929722765037
$rawData.Range("J21").Value = 5.1529178015613
$rawData.Range("L21").Value = 3524.813103493363
$rawData.Range("M21").Value = 1238.449764534581

# --- Summary sheet: final-period rollup per scenario ---
$summary.Range("C2").Value = 460
$summary.Range("D2").Value = 3747.081408254661
$summary.Range("E2").Value = 7.405522824688475
$summary.Range("F2").Value = 57.95892898370582
$summary.Range("C3").Value = 736
$summary.Range("D3").Value = 5167.488790218209
$summary.Range("E3").Value = 6.875513654645529
$summary.Range("F3").Value = 57.6570267755023
$summary.Range("C4").Value = 452
$summary.Range("D4").Value = 2442.758324069539
$summary.Range("E4").Value = 6.707368506691897
$summary.Range("F4").Value = 46.34522099789339
$summary.Range("C5").Value = 469
$summary.Range("D5").Value = 3524.813103493363
$summary.Range("E5").Value = 6.796704856080432
$summary.Range("F5").Value = 63.41420701930006
